# Update "想去人数" (F column) figures for three events that appear on both
# the "展览" sheet (sheet1) and the "全部类型" sheet (sheet4):
#   CM04动漫游戏博览会        4402 -> 4404
#   宋亭山河·炎国明日方舟同人ONLY  93 -> 94
#   云蒸动漫音乐嘉年华        4965 -> 4968

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F8").Value = 4404
$ws1.Range("F9").Value = 94
$ws1.Range("F10").Value = 4968

$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F9").Value = 4404
$ws4.Range("F10").Value = 94
$ws4.Range("F11").Value = 4968
